$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from " FT_L7203" to " FT_L7205"
$ws.Name = " FT_L7205"

# Column A: program path label (row 2-4)
$ws.Range("A2:A4").Value = "FT\L7\L7205五類資產分類上傳轉檔作業"

# Column B: program name (row 2-4)
$ws.Range("B2:B4").Value = "L7205五類資產分類上傳轉檔作業"

# Columns C & D: program code (row 2-4)
$ws.Range("C2:C4").Value = "L7205"
$ws.Range("D2:D4").Value = "L7205"

# Column E: step labels stay the same per row (年/月份/選擇檔案) -- unchanged

# Column F: step descriptions
$ws.Range("F2").Value = "1.限輸入數字，檢核條件：不可輸入0/V(2)"
$ws.Range("F3").Value = "1.限輸入數字，檢核條件：需介於01至12/V(5)2.預設值為上月，若會計日期為月底日則預設本月"
$ws.Range("F4").Value = "1.以滑鼠點選[選擇檔案]按鈕，選取指定位置、檔名之上傳檔案"

# Column K: document reference
$ws.Range("K2:K4").Value = "7-功能規格書"

# Column L: line
$ws.Range("L2:L4").Value = "L7"

# Column M: based-on document name and version
$ws.Range("M2:M4").Value = "製作依據之需求規格書與版本：PJ201800012_URS_7介接外部系統_V1.94.docx"

# Column Q: test case creation date, updated from 2022/01/17 to 2022/02/24
$ws.Range("Q2").Value = 44616
$ws.Range("Q3").Value = 44616
$ws.Range("Q4").Value = 44616
